$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Report generated for the handback of 9126e01e-5a35-4f76-a9b7-f214281ddc0a
# The source handoff file's "Ready for handoff" row (row 7) on both the
# zh-cn and de-de sheets now has a handback result: the target/handback
# files were recorded, but the handback was flagged as stale ("not the
# latest") - same general shape as the "d8c9a66c..." row already present
# earlier in the sheet, only here it's a version/date mismatch instead of
# a file-name mismatch.
# ---------------------------------------------------------------------------

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19bdb02e997071d1c9cbf943503b7a982f0cf786/e2e/9126e01e-5a35-4f76-a9b7-f214281ddc0a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d166b40df505823a6478754ab709267d2321a4d/e2e/9126e01e-5a35-4f76-a9b7-f214281ddc0a.md."

function Set-HandbackRow($ws, $targetFileName, $handbackDateTime, $hyperlinkTarget, $errText) {
    $row = 7

    # Latest Target File (I) - same source markdown file name, rendered as
    # a hyperlink the same way column A / other "Target File" cells are.
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value = $targetFileName
    $iCell.Font.Underline = 2
    $iCell.Font.Color = 15570276
    $iCell.Font.Name = "Calibri"
    $iCell.Font.Size = 11
    $ws.Hyperlinks.Add($iCell, $hyperlinkTarget, "", "", $targetFileName)

    # Latest Handback File (J) - mirrors the Latest Handoff File (G) value.
    $gCell = $ws.Cells.Item($row, 7)
    $ws.Cells.Item($row, 10).Value = $gCell.Value2

    # Latest Handback DateTime (K)
    $ws.Cells.Item($row, 11).Value = $handbackDateTime

    # Error Detail (P)
    $ws.Cells.Item($row, 16).Value = $errText
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn "9126e01e-5a35-4f76-a9b7-f214281ddc0a.md" "2016-08-21 09:00:44" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7d166b40df505823a6478754ab709267d2321a4d/e2e/9126e01e-5a35-4f76-a9b7-f214281ddc0a.md" $errorDetail

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe "9126e01e-5a35-4f76-a9b7-f214281ddc0a.md" "2016-08-21 09:00:50" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7d166b40df505823a6478754ab709267d2321a4d/e2e/9126e01e-5a35-4f76-a9b7-f214281ddc0a.md" $errorDetail
